# Apply the "Added verification email and modified dataplicity status
# read" edit: drop the obsolete/duplicate unit rows from the status
# table and fix up the Cleantech location code.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Delete rows highest-index-first so earlier indices stay valid as we go.
# (row 1 = header; deleted rows, by original 1-based index:
#  13 = "5kW - Unit 12", 12 = "100kW - Unit 11", 11 = "100kW - Unit 10",
#   7 = "5kW - Unit 5",   6 = "5kW - Unit 4",
#   4 = "5kW - Unit 2",   3 = "5kW - Unit 1 Backup 1", 2 = "5kW - Unit 0")
$rowsToDelete = @(13, 12, 11, 7, 6, 4, 3, 2)
foreach ($idx in $rowsToDelete) {
    $t.Rows.Item($idx).Delete()
}

# Fix the Cleantech_10_100kWh row's location from "SG - CT1" to "SG - CT3".
$d.Content.Find.Execute("SG - CT1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SG - CT3", 2)
